$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column S by copying formats (incl. styles/borders) from column R ---
# Row 3 (thin bottom-border strip), row 4 (year header), rows 5-7 (data), row 8 (total, bottom border)
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)

$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- New column S values ---
$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 13.600365850576139
$ws.Range("S6").Value = 9.2742414863791556
$ws.Range("S7").Value = 17.303523954725925
$ws.Range("S8").Value = 205.5

# --- Updated values in existing columns P, Q, R ---
$ws.Range("P5").Value = 23.111083656771282
$ws.Range("Q5").Value = 24.08077930418019
$ws.Range("R5").Value = 19.336931533747723

$ws.Range("P6").Value = 14.322631450320875
$ws.Range("Q6").Value = 13.073459110725862
$ws.Range("R6").Value = 10.464141365743002

$ws.Range("P7").Value = 23.612622725489956

# --- Update selection to match new active cell ---
$ws.Range("Q15").Select() | Out-Null
